# Regenerate the "K" column (column G) values for the save_data sheet.
# The sheet previously stored a raw strike count in column G ("K"); this
# recalculates/regenerates those values (e.g. std/mean-derived s_vals) and
# writes the refreshed numbers back into G2:G67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values, in row order starting at row 2 (row 1 is the header).
$kValues = @(
    2,1,0,0,1,1,2,0,2,2,
    2,1,2,0,1,0,0,1,1,1,
    1,0,1,0,0,0,2,0,1,2,
    0,0,0,2,2,0,0,0,0,0,
    2,3,1,2,0,1,0,0,0,1,
    1,0,0,1,2,0,3,0,0,2,
    1,2,1,1,0,2
)

$startRow = 2
$col = 7  # column G

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, $col).Value = $kValues[$i]
}
